$d = $word.ActiveDocument

# Update the date line (paragraph 1)
$d.Content.Find.Execute("2024-08-24 Saturday", $true, $false, $false, $false, $false,
                        $true, 1, $false, "2024-08-25 Sunday", 2) | Out-Null

# Update the multiplication-table answer cells (first table in the document).
# Cell text is assigned directly (instead of Find/Replace) so that duplicate
# old/new values across different cells cannot cross-contaminate each other.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "518×6=3108"
$t.Cell(1, 2).Range.Text = "767×2=1534"
$t.Cell(1, 3).Range.Text = "931×6=5586"
$t.Cell(1, 4).Range.Text = "632×7=4424"
$t.Cell(1, 5).Range.Text = "624×3=1872"

$t.Cell(5, 1).Range.Text = "218×3=654"
$t.Cell(5, 2).Range.Text = "702×9=6318"
$t.Cell(5, 3).Range.Text = "826×7=5782"
$t.Cell(5, 4).Range.Text = "278×4=1112"
$t.Cell(5, 5).Range.Text = "476×5=2380"

$t.Cell(10, 1).Range.Text = "932×8=7456"
$t.Cell(10, 2).Range.Text = "723×7=5061"
$t.Cell(10, 3).Range.Text = "209×5=1045"
$t.Cell(10, 4).Range.Text = "461×7=3227"
$t.Cell(10, 5).Range.Text = "549×4=2196"

$t.Cell(15, 1).Range.Text = "829×3=2487"
$t.Cell(15, 2).Range.Text = "495×5=2475"
$t.Cell(15, 3).Range.Text = "482×5=2410"
$t.Cell(15, 4).Range.Text = "139×6=834"
$t.Cell(15, 5).Range.Text = "231×7=1617"

$t.Cell(20, 1).Range.Text = "322×4=1288"
$t.Cell(20, 2).Range.Text = "329×7=2303"
$t.Cell(20, 3).Range.Text = "916×4=3664"
$t.Cell(20, 4).Range.Text = "465×5=2325"
$t.Cell(20, 5).Range.Text = "446×7=3122"
